$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '58.436.50'
$ws.Range('E2').Value = '  -1.85%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.616.66'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '533.55'
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.35'
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.566'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.92'
$ws.Range('E9').Value = '  +6.36%  '
$ws.Range('E10').Value = '  -2.28%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('E12').Value = '  +1.03%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.081.12'
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '58.353.84'
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '20.64'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.608.23'
$ws.Range('E16').Value = '  +1.68%  '
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '333.07'
$ws.Range('E19').Value = '  -2.46%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.09'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('E21').Value = '  -2.32%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '66.23'
$ws.Range('E23').Value = '  -1.62%  '
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.32%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.162'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.06'
$ws.Range('E27').Value = '  -2.42%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.998'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0730'
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.63'
$ws.Range('E30').Value = '  -2.46%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.86'
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '18.82'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '150.36'
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.87'
$ws.Range('E34').Value = '  -2.56%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.850'
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.09'
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.809'
$ws.Range('E37').Value = '  -1.92%  '
$ws.Range('E38').Value = '  -3.51%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.55'
$ws.Range('E39').Value = '  +0.41%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '279.62'
$ws.Range('E40').Value = '  +2.19%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '10.68'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '18.91'
$ws.Range('E44').Value = '  +1.92%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0525'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0932'
$ws.Range('E46').Value = '  -2.16%  '
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.934.84'
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.43'
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '17.81'
$ws.Range('E50').Value = '  -4.06%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '113.28'
$ws.Range('E51').Value = '  +1.23%  '
